$d = $word.ActiveDocument

# 1. Strip the "Note" paragraph style off the lone paragraph in the body so it
#    falls back to the default (Normal) style - removes <w:pPr><w:pStyle/></w:pPr>.
foreach ($p in $d.Paragraphs) {
    if ($p.Style.NameLocal -eq "Note") {
        $p.Style = "Normal"
    }
}

# 2. Fix the "MarginNoteRIght" typo -> "MarginNoteRight" on the custom style
#    used for right-hand margin notes (avoids the duplicate-style-name load
#    that confused the i18n.xsl stylesheet lookup).
try {
    $rightStyle = $d.Styles("MarginNoteRIght")
    $rightStyle.NameLocal = "MarginNoteRight"
} catch {
    # style already renamed / not present - nothing to do
}
